$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.676148056983948
$ws.Range("B1").Value = 1.725333333015442
$ws.Range("C1").Value = 1.707668304443359
$ws.Range("D1").Value = 2.041229963302612
$ws.Range("E1").Value = 2.875803709030151
